$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between rows 16 and 18
$ws.Range("E16").Value = "2203"
$ws.Range("E18").Value = "2201"

# Swap the "Valor Mora" values between rows 16 and 18
$ws.Range("F16").Value = 32707
$ws.Range("F18").Value = 36341
